# Update weekly Fruta / hortaliza prices: the data rows (2-10) for the
# Cebollín - Agro Chillán subset get their date/price/origin details
# re-shuffled across rows (as published on a different day), while the
# descriptive columns (A,B,C,E,F,G,H,I,R) stay identical.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that carry the data being re-ordered across rows 2-10.
$cols = @("D", "J", "K", "L", "M", "N", "O", "P", "Q")

# Snapshot the current ("before") values for every affected row/column
# so the subsequent writes don't clobber values we still need to read.
# NOTE: use .Value2 (not .Value) to read raw values here.
$before = @{}
foreach ($row in 2..10) {
    $rowData = @{}
    foreach ($col in $cols) {
        $rowData[$col] = $ws.Range("$col$row").Value2
    }
    $before[$row] = $rowData
}

# Mapping of target row -> source row (which row's data it now holds).
$mapping = @{
    2  = 10
    3  = 5
    4  = 8
    5  = 3
    6  = 6
    7  = 4
    8  = 2
    9  = 7
    10 = 9
}

foreach ($row in 2..10) {
    $srcRow = $mapping[$row]
    $srcData = $before[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$row").Value = $srcData[$col]
    }
}
